$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the old column-B header text ("Jun_13") before we touch it -
# it needs to slide over to the new column D.
$oldHeader = $ws.Range("B1").Value2

# Insert two brand-new date columns before the existing "history" column C.
# This pushes the old column C (and everything in it) two slots to the
# right, landing it in column E.
$ws.Columns("C:D").Insert()

# The engine's column-insert doesn't preserve the custom-width flag on the
# shifted column, so restore explicit widths (raw width 8.0) on all three
# affected columns.
$ws.Columns("C").ColumnWidth = 7.166666667
$ws.Columns("D").ColumnWidth = 7.166666667
$ws.Columns("E").ColumnWidth = 7.166666667

# New columns C and D have no data for any existing analyst rows yet -
# fill them with the same "UN" placeholder used elsewhere.
$ws.Range("C2:D27").Value = "UN"

# Column headers: two new date columns, plus the old column-B header text
# sliding into column D.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = $oldHeader

# Row 22 (BidaskClub) got a brand-new upgrade note in column B, highlighted
# with the same fill used for the other noteworthy rating cells (column E
# now holds the data/formatting that used to live in column C).
$ws.Range("E18").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "6/16/2018,Upgrades,Buy -> Strong-Buy,"
